$d = $word.ActiveDocument

# --- 1. Register the three new paragraph styles (based on MSCParagraph) ---
$baseStyle = $d.Styles("MSCParagraph")

$styleA = $d.Styles.Add("MSC_Paragraph_A", 1)
$styleA.BaseStyle = $baseStyle

$styleB = $d.Styles.Add("MSC_Paragraph_B", 1)
$styleB.BaseStyle = $baseStyle

$styleC = $d.Styles.Add("MSC_Paragraph_C", 1)
$styleC.BaseStyle = $baseStyle

# --- 2. Walk every table in the document and apply the right style ---
for ($t = 1; $t -le $d.Tables.Count; $t++) {
    $tbl = $d.Tables($t)
    $tblStyleName = $tbl.Style.NameLocal

    if ($tblStyleName -eq "MSC_Text_Table_Horiz") {
        # The verse/translation row is the last row of the table (the one
        # that actually holds the MSC_Paragraph bible-text paragraphs).
        $row = $tbl.Rows($tbl.Rows.Count)

        # Column A (text column): only re-style the paragraphs that were
        # using MSC_Paragraph; leave MSC_Join paragraphs untouched.
        $cellA = $row.Cells(1)
        for ($i = 1; $i -le $cellA.Range.Paragraphs.Count; $i++) {
            $p = $cellA.Range.Paragraphs($i)
            if ($p.Style.NameLocal -eq "MSC_Paragraph") {
                $p.Style = "MSC_Paragraph_A"
            }
        }

        # Column B
        $cellB = $row.Cells(2)
        for ($i = 1; $i -le $cellB.Range.Paragraphs.Count; $i++) {
            $cellB.Range.Paragraphs($i).Style = "MSC_Paragraph_B"
        }

        # Column C
        $cellC = $row.Cells(3)
        for ($i = 1; $i -le $cellC.Range.Paragraphs.Count; $i++) {
            $cellC.Range.Paragraphs($i).Style = "MSC_Paragraph_C"
        }
    }
    elseif ($tblStyleName -eq "MSC_Copyright_Table_Horiz") {
        $row = $tbl.Rows(1)
        for ($c = 2; $c -le $row.Cells.Count; $c++) {
            $cell = $row.Cells($c)
            for ($i = 1; $i -le $cell.Range.Paragraphs.Count; $i++) {
                $cell.Range.Paragraphs($i).Style = "MSCCopyright"
            }
        }
    }
}
